$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header-ish values tweaked
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 ("CON"): B2, D2, E2 values were removed (Lichtwark deleted values),
# C2 updated to a new value
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -6.7848162828058793
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 ("STR"): values tweaked
$ws.Range("B3").Value = -7.741317205820601
$ws.Range("C3").Value = -2.4586860437331035
$ws.Range("D3").Value = -13.611244854621201
$ws.Range("E3").Value = 24.051862618030981

# Selection narrowed to reflect the smaller edited range
$ws.Range("B1:E3").Select()
